# Atualização de bases das ligas, do dia: 07-04-2024 às 22:30
#
# Re-sync of match rows for "Mexico Liga de Expansion":
#  - rows 91/92 swap their match data (id/home/away/score/odds)
#  - rows 186/187 swap their match data (id/home/away/score/odds)
#  - a new match (row 224) is appended at the end of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 91 <-> 92 : swap all columns B..AC (column A keeps the row index)
# ---------------------------------------------------------------------

# New row 91 (was row 92's data)
$ws.Cells.Item(91,2).Value2  = 6924568
$ws.Cells.Item(91,6).Value   = "Atletico Morelia"
$ws.Cells.Item(91,7).Value   = "Atlante"
$ws.Cells.Item(91,8).Value2  = 0
$ws.Cells.Item(91,9).Value2  = 1
$ws.Cells.Item(91,10).Value  = "A"
$ws.Cells.Item(91,11).Value2 = 2.4
$ws.Cells.Item(91,12).Value2 = 3
$ws.Cells.Item(91,13).Value2 = 2.875
$ws.Cells.Item(91,14).Value2 = 2.7
$ws.Cells.Item(91,15).Value2 = 3.1
$ws.Cells.Item(91,16).Value2 = 2.8
$ws.Cells.Item(91,17).Value2 = 0
$ws.Cells.Item(91,18).Value2 = 1.85
$ws.Cells.Item(91,19).Value2 = 1.95
$ws.Cells.Item(91,20).Value2 = 2.25
$ws.Cells.Item(91,21).Value2 = 1.975
$ws.Cells.Item(91,22).Value2 = 1.725
$ws.Cells.Item(91,23).Value2 = -1
$ws.Cells.Item(91,24).Value2 = -1
$ws.Cells.Item(91,25).Value2 = 1.8
$ws.Cells.Item(91,26).Value2 = -1
$ws.Cells.Item(91,27).Value2 = 0.95
$ws.Cells.Item(91,28).Value2 = -1
$ws.Cells.Item(91,29).Value2 = 0.7250000000000001

# New row 92 (was row 91's data)
$ws.Cells.Item(92,2).Value2  = 6924569
$ws.Cells.Item(92,6).Value   = "Venados FC"
$ws.Cells.Item(92,7).Value   = "Dorados"
$ws.Cells.Item(92,8).Value2  = 4
$ws.Cells.Item(92,9).Value2  = 1
$ws.Cells.Item(92,10).Value  = "H"
$ws.Cells.Item(92,11).Value2 = 1.615
$ws.Cells.Item(92,12).Value2 = 4
$ws.Cells.Item(92,13).Value2 = 4.5
$ws.Cells.Item(92,14).Value2 = 1.5
$ws.Cells.Item(92,15).Value2 = 4.75
$ws.Cells.Item(92,16).Value2 = 5.75
$ws.Cells.Item(92,17).Value2 = -1.25
$ws.Cells.Item(92,18).Value2 = 1.925
$ws.Cells.Item(92,19).Value2 = 1.875
$ws.Cells.Item(92,20).Value2 = 3
$ws.Cells.Item(92,21).Value2 = 1.75
$ws.Cells.Item(92,22).Value2 = 1.95
$ws.Cells.Item(92,23).Value2 = 0.5
$ws.Cells.Item(92,24).Value2 = -1
$ws.Cells.Item(92,25).Value2 = -1
$ws.Cells.Item(92,26).Value2 = 0.925
$ws.Cells.Item(92,27).Value2 = -1
$ws.Cells.Item(92,28).Value2 = 0.75
$ws.Cells.Item(92,29).Value2 = -1

# ---------------------------------------------------------------------
# Rows 186 <-> 187 : swap all columns B..AC (column A keeps the row index)
# ---------------------------------------------------------------------

# New row 186 (was row 187's data)
$ws.Cells.Item(186,2).Value2  = 7648958
$ws.Cells.Item(186,6).Value   = "Monterrey U23"
$ws.Cells.Item(186,7).Value   = "Mazatlan FC U23"
$ws.Cells.Item(186,8).Value2  = 4
$ws.Cells.Item(186,9).Value2  = 3
$ws.Cells.Item(186,10).Value  = "H"
$ws.Cells.Item(186,11).Value2 = 2.375
$ws.Cells.Item(186,12).Value2 = 3.1
$ws.Cells.Item(186,13).Value2 = 2.75
$ws.Cells.Item(186,14).Value2 = 2.375
$ws.Cells.Item(186,15).Value2 = 3.4
$ws.Cells.Item(186,16).Value2 = 3
$ws.Cells.Item(186,17).Value2 = -0.25
$ws.Cells.Item(186,18).Value2 = 2
$ws.Cells.Item(186,19).Value2 = 1.8
$ws.Cells.Item(186,20).Value2 = 2.75
$ws.Cells.Item(186,21).Value2 = 1.95
$ws.Cells.Item(186,22).Value2 = 1.85
$ws.Cells.Item(186,23).Value2 = 1.375
$ws.Cells.Item(186,24).Value2 = -1
$ws.Cells.Item(186,25).Value2 = -1
$ws.Cells.Item(186,26).Value2 = 1
$ws.Cells.Item(186,27).Value2 = -1
$ws.Cells.Item(186,28).Value2 = 0.95
$ws.Cells.Item(186,29).Value2 = -1

# New row 187 (was row 186's data)
$ws.Cells.Item(187,2).Value2  = 7648957
$ws.Cells.Item(187,6).Value   = "Unam Pumas U23"
$ws.Cells.Item(187,7).Value   = "Tijuana U23"
$ws.Cells.Item(187,8).Value2  = 2
$ws.Cells.Item(187,9).Value2  = 0
$ws.Cells.Item(187,10).Value  = "H"
$ws.Cells.Item(187,11).Value2 = 1.666
$ws.Cells.Item(187,12).Value2 = 3.5
$ws.Cells.Item(187,13).Value2 = 4.2
$ws.Cells.Item(187,14).Value2 = 1.533
$ws.Cells.Item(187,15).Value2 = 4.333
$ws.Cells.Item(187,16).Value2 = 6
$ws.Cells.Item(187,17).Value2 = -1.25
$ws.Cells.Item(187,18).Value2 = 2.025
$ws.Cells.Item(187,19).Value2 = 1.775
$ws.Cells.Item(187,20).Value2 = 2.75
$ws.Cells.Item(187,21).Value2 = 1.775
$ws.Cells.Item(187,22).Value2 = 2.025
$ws.Cells.Item(187,23).Value2 = 0.5329999999999999
$ws.Cells.Item(187,24).Value2 = -1
$ws.Cells.Item(187,25).Value2 = -1
$ws.Cells.Item(187,26).Value2 = 1.025
$ws.Cells.Item(187,27).Value2 = -1
$ws.Cells.Item(187,28).Value2 = -1
$ws.Cells.Item(187,29).Value2 = 1.025

# ---------------------------------------------------------------------
# New row 224 appended at the bottom of the table
# ---------------------------------------------------------------------

# Clone formatting (styles) of the last existing row so the new row gets
# the same cell styles (bordered/bold id column, date-formatted date column).
$ws.Range("A223:AC223").Copy($ws.Range("A224:AC224"))

$ws.Cells.Item(224,1).Value2  = 222
$ws.Cells.Item(224,2).Value2  = 7893216
$ws.Cells.Item(224,3).Value   = "Mexico Liga de Expansion"
$ws.Cells.Item(224,4).Value   = "Mexico Liga de Expansion"
$ws.Cells.Item(224,5).Value2  = 45388.00347222222
$ws.Cells.Item(224,6).Value   = "Dorados"
$ws.Cells.Item(224,7).Value   = "Tapatio"
$ws.Cells.Item(224,8).Value2  = 1
$ws.Cells.Item(224,9).Value2  = 3
$ws.Cells.Item(224,10).Value  = "A"
$ws.Cells.Item(224,11).Value2 = 2.3
$ws.Cells.Item(224,12).Value2 = 3
$ws.Cells.Item(224,13).Value2 = 2.9
$ws.Cells.Item(224,14).Value2 = 3.3
$ws.Cells.Item(224,15).Value2 = 3.6
$ws.Cells.Item(224,16).Value2 = 2.1
$ws.Cells.Item(224,17).Value2 = 0.25
$ws.Cells.Item(224,18).Value2 = 2
$ws.Cells.Item(224,19).Value2 = 1.8
$ws.Cells.Item(224,20).Value2 = 2.5
$ws.Cells.Item(224,21).Value2 = 1.825
$ws.Cells.Item(224,22).Value2 = 1.975
$ws.Cells.Item(224,23).Value2 = -1
$ws.Cells.Item(224,24).Value2 = -1
$ws.Cells.Item(224,25).Value2 = 1.1
$ws.Cells.Item(224,26).Value2 = -1
$ws.Cells.Item(224,27).Value2 = 0.8
$ws.Cells.Item(224,28).Value2 = 0.825
$ws.Cells.Item(224,29).Value2 = -1
